$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column F ("Result"), shifting it to column G
$ws.Columns.Item(6).Insert()

# Copy the header style from the neighboring header cell (E1) onto the new F1 header
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("F1").Value = "Odds"

# Fill in the new "Odds" values for rows 2-9
$ws.Range("F2").Value = 1.01
$ws.Range("F3").Value = 1.11
$ws.Range("F4").Value = 1.18
$ws.Range("F5").Value = 1.84
$ws.Range("F6").Value = 1.84
$ws.Range("F7").Value = 1.84
$ws.Range("F8").Value = 1.84
$ws.Range("F9").Value = 1.84

$excel.CutCopyMode = $false
